$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "System and smoke "

# Update the "Test data" cell for TC_Vtiger_002 (row 7) with the new scenario data
$newTestData = "*URL:http://localhost:8888/`n*User Name: admin                                                                                        *Password: root                                                                                     *lasrtName : Amazon`nindustry: anyone(from the dropdown)`ntypeDropDown : (anyone)"
$ws.Range("G7").Value = $newTestData

# Move the active selection to E9 (also drops the old topLeftCell scroll anchor)
$ws.Range("E9").Select()
